$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.050393171873581, 1.060059378363251, 1.057316900730347, 1.067763724239627, 1.028321640044225, 1.055426563910119, 1.062787189178987, 1.060052202493248, 1.070470713750412, 1.056925390689991),
    @(1.02, 1.051816026289495, 1.061402787084415, 1.058581220172807, 1.069111091203671, 1.028309013115019, 1.056496770783182, 1.063943612724984, 1.061129189832388, 1.071632606921177, 1.057997117379567),
    @(1.02, 1.052736274164227, 1.062271796679359, 1.059399025453466, 1.069982283623907, 1.028298401527434, 1.057188395812965, 1.064691089179896, 1.061825221848498, 1.072383242891091, 1.058689724596223),
    @(1.02, 1.05312304707518, 1.062637068260428, 1.05974276362489, 1.070348382020788, 1.028293355270579, 1.057478949679894, 1.065005138273636, 1.062117632925963, 1.072698529078393, 1.05898069108297),
    @(1.02, 1.05318798226015, 1.062698395547571, 1.059800474930364, 1.070409842799529, 1.028292473646965, 1.057527723003126, 1.065057857447969, 1.062166718342682, 1.072751450612343, 1.059029533669916),
    @(1.02, 1.052741442627236, 1.062276677687691, 1.05940361875972, 1.069987176037478, 1.02829833639928, 1.057192279010843, 1.064695286261965, 1.061829129844997, 1.07238745686184, 1.058693613308687),
    @(1.02, 1.050874124850312, 1.060513446782984, 1.057744246199276, 1.068219208571906, 1.028317877725795, 1.055788427406305, 1.063178176916896, 1.060416353324002, 1.070863628433186, 1.057287768073848),
    @(1.02, 1.047580140409733, 1.057404214168909, 1.054817824444295, 1.065098742193138, 1.028333647786677, 1.053307857150254, 1.060498513964501, 1.057920213741252, 1.068169229251141, 1.054803675123507),
    @(1.02, 1.045381484283809, 1.055329691244991, 1.052865062256223, 1.063014803523231, 1.028331654100369, 1.051649373303894, 1.058707619516667, 1.056251466699201, 1.066366579060366, 1.053142836039855),
    @(1.02, 1.044428743966888, 1.05443094543819, 1.052019018202008, 1.062111526736915, 1.028327832734958, 1.050930055881235, 1.057931040787198, 1.055527735059033, 1.065584455137184, 1.052422497103949),
    @(1.02, 1.044074741748224, 1.054097037079862, 1.05170468281677, 1.061775867424598, 1.028325969513585, 1.050662687744561, 1.057642414140389, 1.055258731953969, 1.065293700683018, 1.052154749273853),
    @(1.02, 1.044150681588195, 1.054168664983895, 1.051772112371805, 1.061847873984978, 1.028326389253746, 1.050720047398793, 1.057704333326612, 1.055316442100172, 1.065356079367465, 1.052212190385374),
    @(1.02, 1.044399484351918, 1.054403346007786, 1.051993036734449, 1.062083783956774, 1.028327687770252, 1.050907958875196, 1.057907186324724, 1.055505502809666, 1.065560426176686, 1.052400368717627),
    @(1.02, 1.044552764927633, 1.054547930771963, 1.052129145206537, 1.062229116920101, 1.028328429040229, 1.051023713154712, 1.058032148007777, 1.055621965793905, 1.065686299233289, 1.052516287381497),
    @(1.02, 1.04544469912281, 1.05538932789263, 1.052921200845854, 1.063074731335027, 1.028331845441816, 1.051697086776844, 1.058759134768262, 1.05629947377354, 1.066418452767081, 1.053190617271409),
    @(1.02, 1.046003992376836, 1.055916986775624, 1.053417903673285, 1.063604914268306, 1.028333196940592, 1.052119156804898, 1.059214854448493, 1.05672414529059, 1.066877291831203, 1.053613286687315),
    @(1.02, 1.046330150774803, 1.056224716840697, 1.053707575496972, 1.063914072758073, 1.028333699793847, 1.052365229236496, 1.059480561239947, 1.056971737990602, 1.067144774199595, 1.053859708569991),
    @(1.02, 1.046441351015917, 1.056329637438464, 1.053806338332634, 1.064019472973071, 1.028333822816739, 1.052449114378395, 1.059571142397187, 1.057056141984758, 1.067235953327209, 1.053943712838413),
    @(1.02, 1.045943992569442, 1.055860378601127, 1.053364617004461, 1.063548039816941, 1.028333081458954, 1.05207388449333, 1.059165971115692, 1.056678593566269, 1.066828078357647, 1.053567950083871),
    @(1.02, 1.044326221278294, 1.054334240378546, 1.051927982173119, 1.062014318306689, 1.028327317635516, 1.050852628673874, 1.057847455918155, 1.055449834057449, 1.065500257757075, 1.052344959941073),
    @(1.02, 1.043308413151836, 1.053374265461935, 1.051024265573987, 1.061049182479762, 1.028321126321882, 1.050083723938094, 1.057017463283442, 1.054676238558637, 1.064664020712006, 1.051574963272373),
    @(1.02, 1.043848036130977, 1.053883208709085, 1.051503386600829, 1.061560898792684, 1.028324651606806, 1.050491435977555, 1.057457553253286, 1.055086434628003, 1.065107458195384, 1.051983254309695),
    @(1.02, 1.045971104109164, 1.055885957535571, 1.05338869509998, 1.063573739205836, 1.028333134522248, 1.05209434145384, 1.059188059719365, 1.056699176772548, 1.066850316273412, 1.053588436095612),
    @(1.02, 1.048432162858885, 1.058208308544898, 1.055574679383319, 1.065906081408108, 1.028331778305837, 1.053949970427282, 1.055574679383319, 1.058566332318337, 1.068866907551088, 1.05544670027504)
)

for ($ri = 0; $ri -lt $data.Count; $ri++) {
    $row = $data[$ri]
    for ($ci = 0; $ci -lt $cols.Count; $ci++) {
        $ws.Cells.Item($ri + 2, $cols[$ci]).Value = $row[$ci]
    }
}
